# ---------------------------------------------------------------------------
# Re-apply the two changes captured in the source commit:
#
#  1. The table on slide 5 gets a different built-in table style
#     (tableStyleId {5A8EFEE3-AFD7-4136-8369-C00D17920B8C} ->
#                   {AA40C37C-7256-4610-93DB-D82839495D1A}).
#
#  2. The deck's applied Design/Theme changes from "Integral" (Red Violet
#     colors) to the default "Office Theme" colors. The two underlying
#     theme parts keep their slots; only the palette that lives in the
#     part referenced by the slide master (the one exposed through the
#     PowerPoint object model as the presentation's ThemeColorScheme)
#     needs to end up holding the Office Theme palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{AA40C37C-7256-4610-93DB-D82839495D1A}")
    }
}

# --- 2. Theme colors: Integral (Red Violet) -> Office Theme ---------------
# Order exposed by ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
